# Applies numeric value updates to the Leve profit-tracking tables
# across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 644.25
$ws.Range("I8").Value = 446.65
$ws.Range("J8").Value = 1632.25
$ws.Range("K8").Value = 1339.95
$ws.Range("L8").Value = 4896.75
$ws.Range("M8").Value = -1200.95
$ws.Range("N8").Value = -5174.75

# Row 116
$ws.Range("H116").Value = 5956.4707
$ws.Range("I116").Value = 5965.3335
$ws.Range("J116").Value = 5935.2
$ws.Range("K116").Value = 5965.3335
$ws.Range("L116").Value = 5935.2
$ws.Range("M116").Value = -2523.3335
$ws.Range("N116").Value = -12819.2

# Row 129
$ws.Range("H129").Value = 116908.336
$ws.Range("I129").Value = 116908.336
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 350725.008
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -345725.008
$ws.Range("N129").Value = $null

# Row 131
$ws.Range("H131").Value = 835499
$ws.Range("I131").Value = 835499
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2506497
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -2501457
$ws.Range("N131").Value = $null

# Row 137
$ws.Range("H137").Value = 12866.37
$ws.Range("J137").Value = 26625
$ws.Range("L137").Value = 79875
$ws.Range("N137").Value = -84975

# Row 138
$ws.Range("H138").Value = 3085.6458
$ws.Range("I138").Value = 3414
$ws.Range("K138").Value = 10242
$ws.Range("M138").Value = -5102

# Row 141
$ws.Range("H141").Value = 4274.5
$ws.Range("I141").Value = 4232.85
$ws.Range("J141").Value = 4378.625
$ws.Range("K141").Value = 12698.55
$ws.Range("L141").Value = 13135.875
$ws.Range("M141").Value = -7518.550000000001
$ws.Range("N141").Value = -23495.875

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 10000000
$ws.Range("I6").Value = 10000000
$ws.Range("K6").Value = 10000000
$ws.Range("M6").Value = -9999827

# Row 32
$ws.Range("H32").Value = 12919.704
$ws.Range("I32").Value = 12753.363
$ws.Range("K32").Value = 12753.363
$ws.Range("M32").Value = -12466.363

# Row 45
$ws.Range("H45").Value = 6556.8237
$ws.Range("I45").Value = 6205.8335
$ws.Range("K45").Value = 6205.8335
$ws.Range("M45").Value = -5828.8335

# Row 61
$ws.Range("H61").Value = 43480.12
$ws.Range("I61").Value = 2489.7273
$ws.Range("K61").Value = 2489.7273
$ws.Range("M61").Value = -2277.7273

# Row 74
$ws.Range("H74").Value = 10067.875
$ws.Range("I74").Value = 1766.5128
$ws.Range("J74").Value = 46040.445
$ws.Range("K74").Value = 1766.5128
$ws.Range("L74").Value = 46040.445
$ws.Range("M74").Value = -892.5128
$ws.Range("N74").Value = -47788.445

# Row 77
$ws.Range("H77").Value = 10067.875
$ws.Range("I77").Value = 1766.5128
$ws.Range("J77").Value = 46040.445
$ws.Range("K77").Value = 8832.564
$ws.Range("L77").Value = 230202.225
$ws.Range("M77").Value = -4464.564
$ws.Range("N77").Value = -238938.225

# Row 136
$ws.Range("H136").Value = 43480.12
$ws.Range("I136").Value = 2489.7273
$ws.Range("K136").Value = 7469.1819
$ws.Range("M136").Value = -4919.1819

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1800.5
$ws.Range("J80").Value = 2153.3845
$ws.Range("L80").Value = 2153.3845
$ws.Range("N80").Value = -4149.3845

# Row 83
$ws.Range("H83").Value = 1800.5
$ws.Range("J83").Value = 2153.3845
$ws.Range("L83").Value = 10766.9225
$ws.Range("N83").Value = -20750.9225

$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null

# Row 31
$ws.Range("H31").Value = 55137.29
$ws.Range("I31").Value = 88845
$ws.Range("J31").Value = 16614.191
$ws.Range("K31").Value = 88845
$ws.Range("L31").Value = 16614.191
$ws.Range("M31").Value = -88550
$ws.Range("N31").Value = -17204.191

# Row 34
$ws.Range("H34").Value = 55137.29
$ws.Range("I34").Value = 88845
$ws.Range("J34").Value = 16614.191
$ws.Range("K34").Value = 88845
$ws.Range("L34").Value = 16614.191
$ws.Range("M34").Value = -88643
$ws.Range("N34").Value = -17018.191

# Row 41
$ws.Range("H41").Value = 5236.3335
$ws.Range("I41").Value = 5236.3335
$ws.Range("K41").Value = 5236.3335
$ws.Range("M41").Value = -4808.3335

# Row 58
$ws.Range("H58").Value = 13158.077
$ws.Range("I58").Value = 3602.1738
$ws.Range("K58").Value = 3602.1738
$ws.Range("M58").Value = -3399.1738

# Row 62
$ws.Range("H62").Value = 4555.2
$ws.Range("J62").Value = 4750
$ws.Range("L62").Value = 4750
$ws.Range("N62").Value = -5998

# Row 65
$ws.Range("H65").Value = 4555.2
$ws.Range("J65").Value = 4750
$ws.Range("L65").Value = 23750
$ws.Range("N65").Value = -29990

# Row 136
$ws.Range("H136").Value = 13158.077
$ws.Range("I136").Value = 3602.1738
$ws.Range("K136").Value = 10806.5214
$ws.Range("M136").Value = -8256.5214

# Row 141
$ws.Range("H141").Value = 73722.25
$ws.Range("J141").Value = 73722.25
$ws.Range("L141").Value = 73722.25
$ws.Range("N141").Value = -84082.25

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5103214
$ws.Range("I5").Value = 1289.625
$ws.Range("J5").Value = 11905780
$ws.Range("K5").Value = 3868.875
$ws.Range("L5").Value = 35717340
$ws.Range("M5").Value = -3756.875
$ws.Range("N5").Value = -35717564

# Row 12
$ws.Range("H12").Value = 62.208332
$ws.Range("I12").Value = 64.23077000000001
$ws.Range("J12").Value = 59.81818
$ws.Range("K12").Value = 192.69231
$ws.Range("L12").Value = 179.45454
$ws.Range("M12").Value = -19.69231000000002
$ws.Range("N12").Value = -525.45454

# Row 13
$ws.Range("H13").Value = 424
$ws.Range("I13").Value = 102.42857
$ws.Range("J13").Value = 1549.5
$ws.Range("K13").Value = 307.28571
$ws.Range("L13").Value = 4648.5
$ws.Range("M13").Value = -139.28571
$ws.Range("N13").Value = -4984.5

# Row 29
$ws.Range("H29").Value = 3685.6667
$ws.Range("I29").Value = 5059.5
$ws.Range("J29").Value = 2998.75
$ws.Range("K29").Value = 15178.5
$ws.Range("L29").Value = 8996.25
$ws.Range("M29").Value = -14901.5
$ws.Range("N29").Value = -9550.25

# Row 33
$ws.Range("H33").Value = 138.33333
$ws.Range("I33").Value = 154
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 924
$ws.Range("L33").Value = 360
$ws.Range("M33").Value = -641
$ws.Range("N33").Value = -926

# Row 68
$ws.Range("H68").Value = 5447.357
$ws.Range("I68").Value = 916
$ws.Range("K68").Value = 2748
$ws.Range("M68").Value = -1937

# Row 71
$ws.Range("H71").Value = 5447.357
$ws.Range("I71").Value = 916
$ws.Range("K71").Value = 8244
$ws.Range("M71").Value = -4188

# Row 129
$ws.Range("H129").Value = 4547493
$ws.Range("J129").Value = 8266788.5
$ws.Range("L129").Value = 24800365.5
$ws.Range("N129").Value = -24810365.5

# Row 135
$ws.Range("H135").Value = 5103214
$ws.Range("I135").Value = 1289.625
$ws.Range("J135").Value = 11905780
$ws.Range("K135").Value = 11606.625
$ws.Range("L135").Value = 107152020
$ws.Range("M135").Value = -9071.625
$ws.Range("N135").Value = -107157090

$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 5580.125
$ws.Range("I99").Value = 2106.8333
$ws.Range("J99").Value = 16000
$ws.Range("K99").Value = 2106.8333
$ws.Range("L99").Value = 16000
$ws.Range("M99").Value = 139.1667000000002
$ws.Range("N99").Value = -20492

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2988407.5
$ws.Range("J7").Value = 12179
$ws.Range("L7").Value = 12179
$ws.Range("N7").Value = -12403

# Row 126
$ws.Range("H126").Value = 2988407.5
$ws.Range("J126").Value = 12179
$ws.Range("L126").Value = 36537
$ws.Range("N126").Value = -41477

# Row 136
$ws.Range("H136").Value = 38875.527
$ws.Range("I136").Value = 41853.12
$ws.Range("J136").Value = 22333.334
$ws.Range("K136").Value = 125559.36
$ws.Range("L136").Value = 67000.00199999999
$ws.Range("M136").Value = -123009.36
$ws.Range("N136").Value = -72100.00199999999

# Row 140
$ws.Range("H140").Value = 147886.89
$ws.Range("J140").Value = 147886.89
$ws.Range("L140").Value = 147886.89
$ws.Range("N140").Value = -158246.89
